$d = $word.ActiveDocument

$replacements = @(
    @("402÷7=57, 3", "210÷9=23, 3"),
    @("897÷2=448, 1", "799÷7=114, 1"),
    @("673÷4=168, 1", "621÷5=124, 1"),
    @("279÷4=69, 3", "664÷3=221, 1"),
    @("952÷6=158, 4", "795÷5=159, 0"),
    @("776÷4=194, 0", "120÷3=40, 0"),
    @("501÷7=71, 4", "802÷7=114, 4"),
    @("308÷2=154, 0", "502÷5=100, 2"),
    @("182÷8=22, 6", "950÷6=158, 2"),
    @("623÷7=89, 0", "557÷5=111, 2"),
    @("123÷6=20, 3", "135÷9=15, 0"),
    @("770÷2=385, 0", "932÷9=103, 5"),
    @("638÷7=91, 1", "881÷4=220, 1"),
    @("876÷7=125, 1", "591÷4=147, 3"),
    @("778÷9=86, 4", "394÷7=56, 2"),
    @("435÷4=108, 3", "271÷5=54, 1"),
    @("250÷2=125, 0", "853÷6=142, 1"),
    @("430÷9=47, 7", "652÷8=81, 4"),
    @("338÷9=37, 5", "167÷7=23, 6"),
    @("115÷6=19, 1", "548÷4=137, 0"),
    @("649÷4=162, 1", "327÷2=163, 1"),
    @("259÷7=37, 0", "280÷4=70, 0"),
    @("210÷4=52, 2", "949÷6=158, 1"),
    @("165÷4=41, 1", "116÷9=12, 8"),
    @("405÷7=57, 6", "357÷5=71, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
